# Updated cryptos list with GitHub Actions: refresh Price (D) and
# Volume(1h) (E) columns for every coin row, and fix the Cronos/Quant
# row ordering (rows 49-50 swap place + data).
# D-column values that look like plain decimal numbers are written with a
# leading "'" so Excel stores/keeps them as text (matching the source
# data, which sometimes has multi-dot thousands separators such as
# "23.414.80" that would never parse as numbers, and must render
# identically to those, e.g. "17.40" not "17.4").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.414.80'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.643.43'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = '''300.27'
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("D7").Value = '''0.3788'
$ws.Range("E7").Value = '  -1.34%  '
$ws.Range("D8").Value = '''50.53'
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("D9").Value = '''0.3507'
$ws.Range("E9").Value = '  -2.63%  '
$ws.Range("D10").Value = '''0.08058'
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("D11").Value = '''1.215'
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("D12").Value = '''1.003'
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").Value = '''22.06'
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").Value = '''6.293'
$ws.Range("E14").Value = '  -2.94%  '
$ws.Range("D15").Value = '''7.244'
$ws.Range("E15").Value = '  -3.38%  '
$ws.Range("D16").Value = '''0.00001209'
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").Value = '1.648.22'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '''95.14'
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '''6.625'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").Value = '''17.40'
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''12.41'
$ws.Range("E23").Value = '  -2.04%  '
$ws.Range("D24").Value = '23.428.42'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").Value = '''2.412'
$ws.Range("E25").Value = '  -4.09%  '
$ws.Range("D26").Value = '''2.976'
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("D27").Value = '''21.04'
$ws.Range("D28").Value = '''151.91'
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("D29").Value = '''5.183'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").Value = '''131.74'
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").Value = '1.830.55'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = '''6.819'
$ws.Range("E32").Value = '  -4.24%  '
$ws.Range("E33").Value = '  -4.77%  '
$ws.Range("E34").Value = '  -7.52%  '
$ws.Range("D35").Value = '''0.9868'
$ws.Range("E35").Value = '  -7.05%  '
$ws.Range("D36").Value = '''0.02684'
$ws.Range("E36").Value = '  -4.03%  '
$ws.Range("D37").Value = '''0.08786'
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").Value = '''5.895'
$ws.Range("E38").Value = '  -3.43%  '
$ws.Range("D39").Value = '''0.2414'
$ws.Range("E39").Value = '  -3.41%  '
$ws.Range("D40").Value = '''0.06780'
$ws.Range("E40").Value = '  -3.17%  '
$ws.Range("D41").Value = '''12.86'
$ws.Range("E41").Value = '  -2.14%  '
$ws.Range("D42").Value = '''0.6867'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").Value = '''1.292'
$ws.Range("E43").Value = '  -3.43%  '
$ws.Range("D44").Value = '''15.55'
$ws.Range("E44").Value = '  -2.67%  '
$ws.Range("D45").Value = '''1.000'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").Value = '''0.6375'
$ws.Range("E46").Value = '  -2.41%  '
$ws.Range("D47").Value = '''3.924'
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("D48").Value = '''2.237'
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '''127.12'
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.07672'
$ws.Range("E50").Value = '  -2.87%  '
$ws.Range("D51").Value = '''1.234'
$ws.Range("E51").Value = '  +2.30%  '
